$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.990.94'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.858.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.13%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5138'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3829'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.73%  '
$ws.Range("E9").Value = '  -9.81%  '
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("E11").Value = '  -0.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.191'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.54'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.861.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.257'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001096'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.59%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.20%  '
$ws.Range("E21").Value = '  +0.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.008'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.020.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.61%  '
$ws.Range("E24").Value = '  -3.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.265'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.074.06'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.508'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.57'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1067'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.24%  '
$ws.Range("E32").Value = '  -3.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.902'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.592'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.400'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02413'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06502'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2179'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6544'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.197'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.993'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.212'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.21%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.16'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6147'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.76%  '
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.671'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.213'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '120.78'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.06'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.36%  '
